$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 205.375  # H19: 244.35 -> 205.375
$ws.Cells.Item(19, 9).Value = 207.25  # I19: 278.0909 -> 207.25
$ws.Cells.Item(19, 10).Value = 203.5  # J19: 203.11111 -> 203.5
$ws.Cells.Item(19, 11).Value = 207.25  # K19: 278.0909 -> 207.25
$ws.Cells.Item(19, 12).Value = 203.5  # L19: 203.11111 -> 203.5
$ws.Cells.Item(19, 13).Value = -32.25  # M19: -103.0909 -> -32.25
$ws.Cells.Item(19, 14).Value = -553.5  # N19: -553.1111100000001 -> -553.5

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 3699.6924  # H43: 1868.25 -> 3699.6924
$ws.Cells.Item(43, 9).Value = 2157  # I43: 1199.5 -> 2157
$ws.Cells.Item(43, 10).Value = 5499.5  # J43: 2002 -> 5499.5
$ws.Cells.Item(43, 11).Value = 2157  # K43: 1199.5 -> 2157
$ws.Cells.Item(43, 12).Value = 5499.5  # L43: 2002 -> 5499.5
$ws.Cells.Item(43, 13).Value = -2088  # M43: -1130.5 -> -2088
$ws.Cells.Item(43, 14).Value = -5637.5  # N43: -2140 -> -5637.5

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1510.6154  # H132: 2978.4 -> 1510.6154
$ws.Cells.Item(132, 9).Value = 1454.3611  # I132: 2906.3572 -> 1454.3611
$ws.Cells.Item(132, 10).Value = 2185.6667  # J132: 3987 -> 2185.6667
$ws.Cells.Item(132, 11).Value = 4363.0833  # K132: 8719.071599999999 -> 4363.0833
$ws.Cells.Item(132, 12).Value = 6557.000100000001  # L132: 11961 -> 6557.000100000001
$ws.Cells.Item(132, 13).Value = -1833.0833  # M132: -6189.071599999999 -> -1833.0833
$ws.Cells.Item(132, 14).Value = -11617.0001  # N132: -17021 -> -11617.0001

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 5282.5454  # H135: 4556.5 -> 5282.5454
$ws.Cells.Item(135, 9).Value = 4545.3335  # I135: 2767.7693 -> 4545.3335
$ws.Cells.Item(135, 10).Value = 6862.2856  # J135: 9207.200000000001 -> 6862.2856
$ws.Cells.Item(135, 11).Value = 40908.0015  # K135: 24909.9237 -> 40908.0015
$ws.Cells.Item(135, 12).Value = 61760.5704  # L135: 82864.8 -> 61760.5704
$ws.Cells.Item(135, 13).Value = -38373.0015  # M135: -22374.9237 -> -38373.0015
$ws.Cells.Item(135, 14).Value = -66830.5704  # N135: -87934.8 -> -66830.5704

# ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 72995.234  # H140: 73859.09 -> 72995.234
$ws.Cells.Item(140, 10).Value = 94193.336  # J140: 94056.25 -> 94193.336
$ws.Cells.Item(140, 12).Value = 94193.336  # L140: 94056.25 -> 94193.336
$ws.Cells.Item(140, 14).Value = -104553.336  # N140: -104416.25 -> -104553.336

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 219926.62  # H32: 4562.86 -> 219926.62
$ws.Cells.Item(32, 9).Value = 8210.25  # I32: 3995.253 -> 8210.25
$ws.Cells.Item(32, 10).Value = 855075.75  # J32: 8361.462 -> 855075.75
$ws.Cells.Item(32, 11).Value = 8210.25  # K32: 3995.253 -> 8210.25
$ws.Cells.Item(32, 12).Value = 855075.75  # L32: 8361.462 -> 855075.75
$ws.Cells.Item(32, 13).Value = -7923.25  # M32: -3708.253 -> -7923.25
$ws.Cells.Item(32, 14).Value = -855649.75  # N32: -8935.462 -> -855649.75

# ARM row 38
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(38, 8).Value = 2000  # H38: 2666.6667 -> 2000

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3354.96  # H61: 4384.1177 -> 3354.96
$ws.Cells.Item(61, 9).Value = 3447.0557  # I61: 4991.1816 -> 3447.0557
$ws.Cells.Item(61, 10).Value = 3118.1428  # J61: 3271.1667 -> 3118.1428
$ws.Cells.Item(61, 11).Value = 3447.0557  # K61: 4991.1816 -> 3447.0557
$ws.Cells.Item(61, 12).Value = 3118.1428  # L61: 3271.1667 -> 3118.1428
$ws.Cells.Item(61, 13).Value = -3235.0557  # M61: -4779.1816 -> -3235.0557
$ws.Cells.Item(61, 14).Value = -3542.1428  # N61: -3695.1667 -> -3542.1428

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 888.5833  # H74: 931.9091 -> 888.5833
$ws.Cells.Item(74, 9).Value = 651.1818  # I74: 675.1 -> 651.1818
$ws.Cells.Item(74, 11).Value = 651.1818  # K74: 675.1 -> 651.1818
$ws.Cells.Item(74, 13).Value = 222.8182  # M74: 198.9 -> 222.8182

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 888.5833  # H77: 931.9091 -> 888.5833
$ws.Cells.Item(77, 9).Value = 651.1818  # I77: 675.1 -> 651.1818
$ws.Cells.Item(77, 11).Value = 3255.909  # K77: 3375.5 -> 3255.909
$ws.Cells.Item(77, 13).Value = 1112.091  # M77: 992.5 -> 1112.091

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3354.96  # H136: 4384.1177 -> 3354.96
$ws.Cells.Item(136, 9).Value = 3447.0557  # I136: 4991.1816 -> 3447.0557
$ws.Cells.Item(136, 10).Value = 3118.1428  # J136: 3271.1667 -> 3118.1428
$ws.Cells.Item(136, 11).Value = 10341.1671  # K136: 14973.5448 -> 10341.1671
$ws.Cells.Item(136, 12).Value = 9354.428400000001  # L136: 9813.500100000001 -> 9354.428400000001
$ws.Cells.Item(136, 13).Value = -7791.167099999999  # M136: -12423.5448 -> -7791.167099999999
$ws.Cells.Item(136, 14).Value = -14454.4284  # N136: -14913.5001 -> -14454.4284

# ARM row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(141, 8).Value = 61084.617  # H141: 61285.715 -> 61084.617
$ws.Cells.Item(141, 10).Value = 64341.668  # J141: 64307.69 -> 64341.668
$ws.Cells.Item(141, 12).Value = 64341.668  # L141: 64307.69 -> 64341.668
$ws.Cells.Item(141, 14).Value = -74701.66800000001  # N141: -74667.69 -> -74701.66800000001

# BSM row 19
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 0  # H19: 6995 -> 0
$ws.Cells.Item(19, 10).Value = 0  # J19: 6995 -> 0
$ws.Cells.Item(19, 12).Value = 0  # L19: 6995 -> 0
$ws.Cells.Item(19, 14).ClearContents()  # N19: -7341 -> (blank)

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1179  # H94: 1171.8857 -> 1179
$ws.Cells.Item(94, 9).Value = 1121.069  # I94: 1120.3667 -> 1121.069
$ws.Cells.Item(94, 10).Value = 1515  # J94: 1481 -> 1515
$ws.Cells.Item(94, 11).Value = 1121.069  # K94: 1120.3667 -> 1121.069
$ws.Cells.Item(94, 12).Value = 1515  # L94: 1481 -> 1515
$ws.Cells.Item(94, 13).Value = -670.069  # M94: -669.3667 -> -670.069
$ws.Cells.Item(94, 14).Value = -2417  # N94: -2383 -> -2417

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6614.32  # H134: 6893.5835 -> 6614.32
$ws.Cells.Item(134, 9).Value = 893.8261  # I134: 940.2857 -> 893.8261
$ws.Cells.Item(134, 10).Value = 72400  # J134: 48566.668 -> 72400
$ws.Cells.Item(134, 11).Value = 2681.4783  # K134: 2820.8571 -> 2681.4783
$ws.Cells.Item(134, 12).Value = 217200  # L134: 145700.004 -> 217200
$ws.Cells.Item(134, 13).Value = -146.4782999999998  # M134: -285.8571000000002 -> -146.4782999999998
$ws.Cells.Item(134, 14).Value = -222270  # N134: -150770.004 -> -222270

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 89300  # H140: 89433.336 -> 89300
$ws.Cells.Item(140, 10).Value = 89300  # J140: 89433.336 -> 89300
$ws.Cells.Item(140, 12).Value = 89300  # L140: 89433.336 -> 89300
$ws.Cells.Item(140, 14).Value = -99660  # N140: -99793.336 -> -99660

# CRP row 39
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(39, 8).Value = 32230.1  # H39: 3825.3333 -> 32230.1
$ws.Cells.Item(39, 9).Value = 1450.25  # I39: 2504.6667 -> 1450.25
$ws.Cells.Item(39, 10).Value = 52750  # J39: 6466.6665 -> 52750
$ws.Cells.Item(39, 11).Value = 1450.25  # K39: 2504.6667 -> 1450.25
$ws.Cells.Item(39, 12).Value = 52750  # L39: 6466.6665 -> 52750
$ws.Cells.Item(39, 13).Value = -1059.25  # M39: -2113.6667 -> -1059.25
$ws.Cells.Item(39, 14).Value = -53532  # N39: -7248.6665 -> -53532

# CRP row 49
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(49, 8).Value = 32230.1  # H49: 3825.3333 -> 32230.1
$ws.Cells.Item(49, 9).Value = 1450.25  # I49: 2504.6667 -> 1450.25
$ws.Cells.Item(49, 10).Value = 52750  # J49: 6466.6665 -> 52750
$ws.Cells.Item(49, 11).Value = 1450.25  # K49: 2504.6667 -> 1450.25
$ws.Cells.Item(49, 12).Value = 52750  # L49: 6466.6665 -> 52750
$ws.Cells.Item(49, 13).Value = -1268.25  # M49: -2322.6667 -> -1268.25
$ws.Cells.Item(49, 14).Value = -53114  # N49: -6830.6665 -> -53114

# CRP row 138
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(138, 8).Value = 49812.5  # H138: 49885.715 -> 49812.5
$ws.Cells.Item(138, 10).Value = 49812.5  # J138: 49885.715 -> 49812.5
$ws.Cells.Item(138, 12).Value = 49812.5  # L138: 49885.715 -> 49812.5
$ws.Cells.Item(138, 14).Value = -60092.5  # N138: -60165.715 -> -60092.5

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(44, 8).Value = 2994.353  # H44: 445.7143 -> 2994.353
$ws.Cells.Item(44, 9).Value = 250  # I44: 353.33334 -> 250
$ws.Cells.Item(44, 10).Value = 3360.2666  # J44: 1000 -> 3360.2666
$ws.Cells.Item(44, 11).Value = 750  # K44: 1060.00002 -> 750
$ws.Cells.Item(44, 12).Value = 10080.7998  # L44: 3000 -> 10080.7998
$ws.Cells.Item(44, 13).Value = -352  # M44: -662.0000199999999 -> -352
$ws.Cells.Item(44, 14).Value = -10876.7998  # N44: -3796 -> -10876.7998

# CUL row 69
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 3410.9092  # H69: 30000 -> 3410.9092
$ws.Cells.Item(69, 10).Value = 3410.9092  # J69: 30000 -> 3410.9092
$ws.Cells.Item(69, 12).Value = 10232.7276  # L69: 90000 -> 10232.7276
$ws.Cells.Item(69, 14).Value = -11854.7276  # N69: -91622 -> -11854.7276

# CUL row 72
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 3410.9092  # H72: 30000 -> 3410.9092
$ws.Cells.Item(72, 10).Value = 3410.9092  # J72: 30000 -> 3410.9092
$ws.Cells.Item(72, 12).Value = 30698.1828  # L72: 270000 -> 30698.1828
$ws.Cells.Item(72, 14).Value = -38810.1828  # N72: -278112 -> -38810.1828

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 463.88  # H107: 529.913 -> 463.88
$ws.Cells.Item(107, 9).Value = 263.75  # I107: 281.42856 -> 263.75
$ws.Cells.Item(107, 10).Value = 558.05884  # J107: 916.44446 -> 558.05884
$ws.Cells.Item(107, 11).Value = 791.25  # K107: 844.28568 -> 791.25
$ws.Cells.Item(107, 12).Value = 1674.17652  # L107: 2749.33338 -> 1674.17652
$ws.Cells.Item(107, 13).Value = 1128.75  # M107: 1075.71432 -> 1128.75
$ws.Cells.Item(107, 14).Value = -5514.17652  # N107: -6589.33338 -> -5514.17652

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 883.57574  # H131: 5377207 -> 883.57574
$ws.Cells.Item(131, 9).Value = 723.7895  # I131: 918.7619 -> 723.7895
$ws.Cells.Item(131, 10).Value = 921.525  # J131: 6945291 -> 921.525
$ws.Cells.Item(131, 11).Value = 2171.3685  # K131: 2756.2857 -> 2171.3685
$ws.Cells.Item(131, 12).Value = 2764.575  # L131: 20835873 -> 2764.575
$ws.Cells.Item(131, 13).Value = 2868.6315  # M131: 2283.7143 -> 2868.6315
$ws.Cells.Item(131, 14).Value = -12844.575  # N131: -20845953 -> -12844.575

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1762.025  # H132: 1984.2632 -> 1762.025
$ws.Cells.Item(132, 9).Value = 1099.4286  # I132: 1278.0625 -> 1099.4286
$ws.Cells.Item(132, 10).Value = 2494.3684  # J132: 2497.8635 -> 2494.3684
$ws.Cells.Item(132, 11).Value = 9894.857399999999  # K132: 11502.5625 -> 9894.857399999999
$ws.Cells.Item(132, 12).Value = 22449.3156  # L132: 22480.7715 -> 22449.3156
$ws.Cells.Item(132, 13).Value = -7364.857399999999  # M132: -8972.5625 -> -7364.857399999999
$ws.Cells.Item(132, 14).Value = -27509.3156  # N132: -27540.7715 -> -27509.3156

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 4435  # H18: 0 -> 4435
$ws.Cells.Item(18, 9).Value = 4152.5  # I18: 0 -> 4152.5
$ws.Cells.Item(18, 10).Value = 5000  # J18: 0 -> 5000
$ws.Cells.Item(18, 11).Value = 4152.5  # K18: 0 -> 4152.5
$ws.Cells.Item(18, 12).Value = 5000  # L18: 0 -> 5000
$ws.Cells.Item(18, 13).Value = -3859.5  # M18: None -> -3859.5
$ws.Cells.Item(18, 14).Value = -5586  # N18: None -> -5586

# GSM row 35
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 0  # H35: 6500 -> 0
$ws.Cells.Item(35, 10).Value = 0  # J35: 6500 -> 0
$ws.Cells.Item(35, 12).Value = 0  # L35: 6500 -> 0
$ws.Cells.Item(35, 14).ClearContents()  # N35: -7096 -> (blank)

# GSM row 58
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(58, 8).Value = 18625  # H58: 17941.176 -> 18625

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3016.3  # H102: 1328.9143 -> 3016.3
$ws.Cells.Item(102, 9).Value = 2872  # I102: 1219.742 -> 2872
$ws.Cells.Item(102, 10).Value = 3593.5  # J102: 2175 -> 3593.5
$ws.Cells.Item(102, 11).Value = 2872  # K102: 1219.742 -> 2872
$ws.Cells.Item(102, 12).Value = 3593.5  # L102: 2175 -> 3593.5
$ws.Cells.Item(102, 13).Value = -1250  # M102: 402.258 -> -1250
$ws.Cells.Item(102, 14).Value = -6837.5  # N102: -5419 -> -6837.5

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(138, 8).Value = 69150  # H138: 69033.336 -> 69150
$ws.Cells.Item(138, 10).Value = 69150  # J138: 69033.336 -> 69150
$ws.Cells.Item(138, 12).Value = 69150  # L138: 69033.336 -> 69150
$ws.Cells.Item(138, 14).Value = -79430  # N138: -79313.336 -> -79430

# GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(140, 8).Value = 89799.5  # H140: 89844.5 -> 89799.5
$ws.Cells.Item(140, 10).Value = 89799.5  # J140: 89844.5 -> 89799.5
$ws.Cells.Item(140, 12).Value = 89799.5  # L140: 89844.5 -> 89799.5
$ws.Cells.Item(140, 14).Value = -100159.5  # N140: -100204.5 -> -100159.5

# LTW row 57
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(57, 8).Value = 243008.2  # H57: 401666.66 -> 243008.2
$ws.Cells.Item(57, 9).Value = 302010.25  # I57: 401666.66 -> 302010.25
$ws.Cells.Item(57, 10).Value = 7000  # J57: 0 -> 7000
$ws.Cells.Item(57, 11).Value = 302010.25  # K57: 401666.66 -> 302010.25
$ws.Cells.Item(57, 12).Value = 7000  # L57: 0 -> 7000
$ws.Cells.Item(57, 13).Value = -301444.25  # M57: -401100.66 -> -301444.25
$ws.Cells.Item(57, 14).Value = -8132  # N57: None -> -8132

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2035.3334  # H68: 2055.276 -> 2035.3334
$ws.Cells.Item(68, 9).Value = 1986.4546  # I68: 1995.3334 -> 1986.4546
$ws.Cells.Item(68, 10).Value = 2133.0908  # J68: 2212.625 -> 2133.0908
$ws.Cells.Item(68, 11).Value = 1986.4546  # K68: 1995.3334 -> 1986.4546
$ws.Cells.Item(68, 12).Value = 2133.0908  # L68: 2212.625 -> 2133.0908
$ws.Cells.Item(68, 13).Value = -1237.4546  # M68: -1246.3334 -> -1237.4546
$ws.Cells.Item(68, 14).Value = -3631.0908  # N68: -3710.625 -> -3631.0908

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 2035.3334  # H71: 2055.276 -> 2035.3334
$ws.Cells.Item(71, 9).Value = 1986.4546  # I71: 1995.3334 -> 1986.4546
$ws.Cells.Item(71, 10).Value = 2133.0908  # J71: 2212.625 -> 2133.0908
$ws.Cells.Item(71, 11).Value = 9932.273000000001  # K71: 9976.666999999999 -> 9932.273000000001
$ws.Cells.Item(71, 12).Value = 10665.454  # L71: 11063.125 -> 10665.454
$ws.Cells.Item(71, 13).Value = -6188.273000000001  # M71: -6232.666999999999 -> -6188.273000000001
$ws.Cells.Item(71, 14).Value = -18153.454  # N71: -18551.125 -> -18153.454

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3060.2693  # H136: 3114.96 -> 3060.2693
$ws.Cells.Item(136, 9).Value = 3030.5881  # I136: 3013.2942 -> 3030.5881
$ws.Cells.Item(136, 10).Value = 3116.3333  # J136: 3331 -> 3116.3333
$ws.Cells.Item(136, 11).Value = 9091.764299999999  # K136: 9039.882599999999 -> 9091.764299999999
$ws.Cells.Item(136, 12).Value = 9348.999899999999  # L136: 9993 -> 9348.999899999999
$ws.Cells.Item(136, 13).Value = -6541.764299999999  # M136: -6489.882599999999 -> -6541.764299999999
$ws.Cells.Item(136, 14).Value = -14448.9999  # N136: -15093 -> -14448.9999

# LTW row 138
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(138, 8).Value = 59817.445  # H138: 58485.7 -> 59817.445
$ws.Cells.Item(138, 10).Value = 59817.445  # J138: 58485.7 -> 59817.445
$ws.Cells.Item(138, 12).Value = 59817.445  # L138: 58485.7 -> 59817.445
$ws.Cells.Item(138, 14).Value = -70097.44500000001  # N138: -68765.7 -> -70097.44500000001

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(139, 8).Value = 56880  # H139: 54216.668 -> 56880
$ws.Cells.Item(139, 10).Value = 69850  # J139: 64060 -> 69850
$ws.Cells.Item(139, 12).Value = 69850  # L139: 64060 -> 69850
$ws.Cells.Item(139, 14).Value = -80130  # N139: -74340 -> -80130

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 9000  # H2: 5563666 -> 9000
$ws.Cells.Item(2, 9).Value = 9000  # I2: 25004990 -> 9000
$ws.Cells.Item(2, 10).Value = 0  # J2: 9002.143 -> 0
$ws.Cells.Item(2, 11).Value = 9000  # K2: 25004990 -> 9000
$ws.Cells.Item(2, 12).Value = 0  # L2: 9002.143 -> 0
$ws.Cells.Item(2, 13).Value = -8888  # M2: -25004878 -> -8888
$ws.Cells.Item(2, 14).ClearContents()  # N2: -9226.143 -> (blank)

# WVR row 29
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 20000000  # H29: 2004088.8 -> 20000000
$ws.Cells.Item(29, 9).Value = 20000000  # I29: 10000400 -> 20000000
$ws.Cells.Item(29, 10).Value = 0  # J29: 5011 -> 0
$ws.Cells.Item(29, 11).Value = 20000000  # K29: 10000400 -> 20000000
$ws.Cells.Item(29, 12).Value = 0  # L29: 5011 -> 0
$ws.Cells.Item(29, 13).Value = -19999710  # M29: -10000110 -> -19999710
$ws.Cells.Item(29, 14).ClearContents()  # N29: -5591 -> (blank)

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 132214.5  # H46: 75685.8 -> 132214.5
$ws.Cells.Item(46, 10).Value = 132214.5  # J46: 75685.8 -> 132214.5
$ws.Cells.Item(46, 12).Value = 132214.5  # L46: 75685.8 -> 132214.5
$ws.Cells.Item(46, 14).Value = -132676.5  # N46: -76147.8 -> -132676.5

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 5684.615  # H54: 6471.32 -> 5684.615
$ws.Cells.Item(54, 9).Value = 1750  # I54: 0 -> 1750
$ws.Cells.Item(54, 10).Value = 6400  # J54: 6471.32 -> 6400
$ws.Cells.Item(54, 11).Value = 1750  # K54: 0 -> 1750
$ws.Cells.Item(54, 12).Value = 6400  # L54: 6471.32 -> 6400
$ws.Cells.Item(54, 13).Value = -1230  # M54: None -> -1230
$ws.Cells.Item(54, 14).Value = -7440  # N54: -7511.32 -> -7440

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(134, 8).Value = 132214.5  # H134: 75685.8 -> 132214.5
$ws.Cells.Item(134, 10).Value = 132214.5  # J134: 75685.8 -> 132214.5
$ws.Cells.Item(134, 12).Value = 396643.5  # L134: 227057.4 -> 396643.5
$ws.Cells.Item(134, 14).Value = -401713.5  # N134: -232127.4 -> -401713.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1086.9736  # H136: 1085.7106 -> 1086.9736
$ws.Cells.Item(136, 9).Value = 850.2273  # I136: 848.0454999999999 -> 850.2273
$ws.Cells.Item(136, 11).Value = 2550.6819  # K136: 2544.1365 -> 2550.6819
$ws.Cells.Item(136, 13).Value = -0.6819000000000415  # M136: 5.863500000000386 -> -0.6819000000000415

# WVR row 139
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(139, 8).Value = 59840  # H139: 58116.668 -> 59840
$ws.Cells.Item(139, 10).Value = 59840  # J139: 58116.668 -> 59840
$ws.Cells.Item(139, 12).Value = 59840  # L139: 58116.668 -> 59840
$ws.Cells.Item(139, 14).Value = -70120  # N139: -68396.66800000001 -> -70120
